$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the expected-result column (C) values: "pass" -> "Pass", "fail" -> "Fail"
$ws.Range("C2").Value = "Pass"
$ws.Range("C3").Value = "Fail"
$ws.Range("C4").Value = "Fail"
$ws.Range("C5").Value = "Fail"

# Update the active selection to C3
$ws.Range("C3").Select()
